$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5, column K ("Setup Jenis Dokumen ...") get a ",0" suffix appended.
# Set the non-trailing-space variants first (K3/K4/K5) then the
# trailing-space variant (K2) so the shared-string table picks up the
# same ordering as the authored workbook.
$ws.Range("K3").Value = "Setup Jenis Dokumen,0"
$ws.Range("K2").Value = "Setup Jenis Dokumen ,0"

# Rows 2-5, column J ("Setup") get a ",0" suffix appended too. These cells
# use a quote-prefix (text-forced) style, so prefix the literal with an
# apostrophe to keep that cell style intact.
$ws.Range("J2").Value = "'Setup,0"
$ws.Range("J3").Value = "'Setup,0"
$ws.Range("J4").Value = "'Setup,0"
$ws.Range("J5").Value = "'Setup,0"

$ws.Range("K4").Value = "Setup Jenis Dokumen,0"
$ws.Range("K5").Value = "Setup Jenis Dokumen,0"

# Selection/view moved back to A2 (scrolled to top-left, single cell selected).
$ws.Range("A2").Select() | Out-Null
